# Remove the leftover "Equipe de Desenvolvimento: " run from the
# (now empty) paragraph that precedes the heading paragraph of the
# same text. Only one run in the document carries this exact text
# (trailing colon + space), so a plain Find/Replace with an empty
# replacement string removes the run's text content without touching
# the sibling heading paragraph ("Equipe de Desenvolvimento", no colon).

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Equipe de Desenvolvimento: ",  # FindText
    $true,                          # MatchCase
    $false,                         # MatchWholeWord
    $false,                         # MatchWildcards
    $false,                         # MatchSoundsLike
    $false,                         # MatchAllWordForms
    $true,                          # Forward
    1,                              # Wrap (wdFindContinue)
    $false,                         # Format
    "",                             # ReplaceWith
    2                               # Replace (wdReplaceAll)
)
